$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) columns; leading apostrophe forces text storage
# so values match the original inline-string cell type instead of becoming numbers.
$ws.Range("D2").Value = "'275.06"
$ws.Range("E2").Value = "'-0.73%"
$ws.Range("E3").Value = "'-2.25%"
$ws.Range("D4").Value = "'4.862"
$ws.Range("E4").Value = "'1.78%"
$ws.Range("D5").Value = "'0.06344"
$ws.Range("E5").Value = "'1.23%"
$ws.Range("D6").Value = "'6.887"
$ws.Range("E6").Value = "'-0.50%"
$ws.Range("D7").Value = "'3.314"
$ws.Range("E7").Value = "'1.46%"
$ws.Range("D8").Value = "'1.253"
$ws.Range("E8").Value = "'32.96%"
$ws.Range("D9").Value = "'0.8677"
$ws.Range("E9").Value = "'-1.32%"
$ws.Range("D10").Value = "'0.1553"
$ws.Range("E10").Value = "'6.38%"
$ws.Range("D11").Value = "'0.05215"
$ws.Range("E11").Value = "'-2.56%"
$ws.Range("D12").Value = "'0.07403"
$ws.Range("E12").Value = "'1.83%"
$ws.Range("D13").Value = "'0.02922"
$ws.Range("E13").Value = "'-5.86%"
$ws.Range("D14").Value = "'0.09031"
$ws.Range("E14").Value = "'-0.36%"
$ws.Range("D15").Value = "'0.001574"
$ws.Range("E15").Value = "'1.30%"
$ws.Range("D16").Value = "'0.0006315"
$ws.Range("E16").Value = "'0.76%"
$ws.Range("D17").Value = "'0.005965"
$ws.Range("E17").Value = "'2.82%"
$ws.Range("E18").Value = "'0.10%"
$ws.Range("D19").Value = "'2.272"
$ws.Range("E19").Value = "'-0.52%"
$ws.Range("D20").Value = "'0.3113"
$ws.Range("E21").Value = "'1.60%"
$ws.Range("D22").Value = "'3.911"
$ws.Range("E22").Value = "'1.61%"
$ws.Range("D23").Value = "'0.04367"
$ws.Range("E23").Value = "'0.83%"
$ws.Range("D24").Value = "'0.001178"
$ws.Range("E24").Value = "'-0.17%"
$ws.Range("D25").Value = "'0.004249"
$ws.Range("E25").Value = "'-0.60%"
$ws.Range("E26").Value = "'-0.22%"
$ws.Range("E27").Value = "'-0.83%"
$ws.Range("D40").Value = "'0.04110"
$ws.Range("E40").Value = "'1.89%"
$ws.Range("D41").Value = "'0.006888"
$ws.Range("E41").Value = "'7.20%"
$ws.Range("D42").Value = "'0.1168"
$ws.Range("E42").Value = "'1.27%"
$ws.Range("E43").Value = "'-2.73%"
$ws.Range("D44").Value = "'0.01075"
$ws.Range("E44").Value = "'-10.92%"
$ws.Range("D45").Value = "'0.00005301"
$ws.Range("E45").Value = "'3.89%"
$ws.Range("D46").Value = "'0.01999"
$ws.Range("E46").Value = "'-33.07%"
$ws.Range("D47").Value = "'1.486"
$ws.Range("E47").Value = "'-37.47%"

# Rows 46 and 47 swapped coin entries (BOLO <-> CoinbaseStockToken); update name/link too
$ws.Range("B46").Value = "CoinbaseStockToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
